# Default title is taken from yt.title
#
# Previously, users typed a free-form (often lowercase) "Title" for each
# download row. Going forward, the default Title is populated from the
# YouTube video's own title (yt.title), which is properly capitalized.
# Normalize the existing lowercase, hand-typed titles in column C to
# Title Case so the sheet is consistent with titles that will now come
# from yt.title. Proper titles (already capitalized) are left untouched.

function ConvertTo-TitleCaseManual {
    param([string]$s)

    if ([string]::IsNullOrEmpty($s)) { return $s }

    $words = $s.Split(" ")
    $result = @()
    foreach ($w in $words) {
        if ($w.Length -gt 0) {
            $first = $w.Substring(0, 1).ToUpper()
            $rest = $w.Substring(1)
            $result += ($first + $rest)
        }
        else {
            $result += $w
        }
    }
    return [string]::Join(" ", $result)
}

function Test-StartsLowercase {
    param([string]$s)

    if ([string]::IsNullOrEmpty($s)) { return $false }
    $code = [int][char]$s[0]
    return ($code -ge 97) -and ($code -le 122)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Title") holds one row per downloaded track, starting at row 2.
# Rewrite every hand-typed, all-lowercase title to Title Case; leave titles
# that are already properly cased (e.g. "Faded", "The Humma Song") alone.
$lastRow = 34
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $title = $cell.Value2
    if (Test-StartsLowercase $title) {
        $cell.Value = ConvertTo-TitleCaseManual $title
    }
}

# Reset the frozen-pane view back to the top of the sheet (it had drifted to
# show row 8 with C9 selected) so the sheet reopens scrolled to the top.
$ws.Range("A2").Select()
